# Update "Vehicle Loan Default Prediction.pptx" title slide:
#   1. Move the existing team-members textbox ("TextBox 6") slightly.
#   2. Add a "Guided by " textbox above it.
#   3. Add a red "Animesh Tiweri" textbox (right aligned) next to "Guided by ".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Reposition the existing team-members shape ("TextBox 6") ---
$teamShape = $s.Shapes.Item("TextBox 6")
$teamShape.Left = 78
$teamShape.Top = 298.51543307086615

# --- Reserve the two lowest free shape ids so the new textboxes created below
#     land on ids 9 / 10 ("TextBox 8" / "TextBox 9"), matching real PowerPoint's
#     "next id after highest used so far" behaviour on this slide (max id is 8).
$reserved1 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$reserved2 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)

# --- 2. Add the "Guided by " textbox ---
$guidedBy = $s.Shapes.AddTextbox(1, 84.12629921259843, 258.6159842519685, 611.8737007874016, 41.198425196850394)
$guidedBy.TextFrame.WordWrap = -1
$guidedBy.TextFrame.AutoSize = 1
$guidedBy.Fill.Visible = 0
$guidedBy.TextFrame.TextRange.Text = "Guided by "
$guidedBy.TextFrame.TextRange.Font.Size = 28

# --- 3. Add the "Animesh Tiweri" textbox (red, right aligned, 3 runs) ---
$guideName = $s.Shapes.AddTextbox(1, 63.67165354330709, 298.57732283464566, 156, 31.50472440944882)
$guideName.TextFrame.WordWrap = -1
$guideName.TextFrame.AutoSize = 1
$guideName.Fill.Visible = 0

$nameRange = $guideName.TextFrame.TextRange
$nameRange.Text = "Animesh"
$nameRange.ParagraphFormat.Alignment = 3
$nameRange.Font.Size = 20
$nameRange.Font.Color.RGB = 255

$spaceRange = $nameRange.InsertAfter(" ")
$spaceRange.Font.Size = 20
$spaceRange.Font.Color.RGB = 255

$lastRange = $spaceRange.InsertAfter("Tiweri")
$lastRange.Font.Size = 20
$lastRange.Font.Color.RGB = 255

# --- Drop the two reservation placeholders now that the ids are claimed ---
$reserved1.Delete()
$reserved2.Delete()
